# Applies the "Statistics section of report" edit:
#  1. Moves the _GoBack bookmark from the XKCD-image paragraph to sit right
#     at the very start of the document (before the H1 run).
#  2. Collapses several runs that were split mid-word (an artifact of how
#     the prose was originally typed/pasted) back into single runs, without
#     changing the visible text.
#  3. Regroups the markdown image-link runs: "!" + "[via XKCD](/01_I" +
#     "mages/...)" -> "![" + "via XKCD](/01_Images/...)" (visible text is
#     unchanged, only where the run boundary -- and so the gramStart
#     proofErr marker -- falls changes).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: move the _GoBack bookmark to the very start of the document.
# ---------------------------------------------------------------------
# Directly adding a collapsed bookmark at absolute position 0 mis-places
# the bookmarkEnd tag, so insert a one-character placeholder first, anchor
# the (now non-zero, unambiguous) bookmark position right after it, then
# delete the placeholder again -- the bookmark slides back to position 0
# cleanly.
$startRng = $d.Range(0, 0)
$startRng.InsertBefore("X")

$bmRng = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRng)

$placeholderRng = $d.Range(0, 1)
$placeholderRng.Delete()

# ---------------------------------------------------------------------
# Step 2: merge the split runs in the two "Initial Ideas" paragraphs.
# ---------------------------------------------------------------------
# Word's Range.Text setter re-writes the covered runs as a single run, so
# re-assigning the same text collapses the multi-run spans into one run.
# Because the replacement text here is identical to what's already there,
# go via a temporary placeholder string so the engine actually performs a
# (merging) write instead of treating it as a no-op.

function Merge-Range-Text($findText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $ok) {
        throw "Merge-Range-Text: text not found: $findText"
    }
    $rng.Text = "TEMP_MERGE_PLACEHOLDER"
    $rng2 = $d.Content
    $ok2 = $rng2.Find.Execute("TEMP_MERGE_PLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $ok2) {
        throw "Merge-Range-Text: placeholder not found for: $findText"
    }
    $rng2.Text = $findText
}

Merge-Range-Text("In starting this project I knew fairly early on the genre (or type) of game that I wanted  to create: a puzzle game")
Merge-Range-Text("This would satisfy the 'educational' criteria for the solution as a game requiring logic and specific thinking processes would be stimulating for all ages, including our target market, which is High School Students.")
Merge-Range-Text("In the discussion amongst the class, we went through a variety of games that we could make")
Merge-Range-Text("Many of these suggestions, however, did not appeal to me as they were merely remakes of existing games")
Merge-Range-Text("Some of those mentioned were Minesweeper, a game bundled with Windows since 1989")

# ---------------------------------------------------------------------
# Step 3: fix the markdown image-link text.
# ---------------------------------------------------------------------
# "!" -> "!["
$rng = $d.Content
$ok = $rng.Find.Execute("!", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $ok) {
    throw "'!' run not found"
}
$rng.Text = "!["

# "[via XKCD](/01_I" + "mages/001_LabyrinthPuzzle_XKCD)" -> "via XKCD](/01_Images/001_LabyrinthPuzzle_XKCD)"
$rng = $d.Content
$ok = $rng.Find.Execute("[via XKCD](/01_Images/001_LabyrinthPuzzle_XKCD)", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $ok) {
    throw "image markdown text not found"
}
$rng.Text = "via XKCD](/01_Images/001_LabyrinthPuzzle_XKCD)"
